# Update the cryptocurrency price/volume table with the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'40.146.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.17%  "

# Row 3
$ws.Range("D3").Value = "'2.209.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.50%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'297.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.18%  "

# Row 6
$ws.Range("D6").Value = "'87.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "

# Row 7
$ws.Range("D7").Value = "'0.514"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.59%  "

# Row 8
$ws.Range("E8").Value = "  -0.16%  "

# Row 9
$ws.Range("D9").Value = "'0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.23%  "

# Row 10
$ws.Range("D10").Value = "'52.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.39%  "

# Row 11
$ws.Range("D11").Value = "'31.01"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.84%  "

# Row 12
$ws.Range("D12").Value = "'0.0782"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.56%  "

# Row 13
$ws.Range("E13").Value = "  +2.00%  "

# Row 14
$ws.Range("D14").Value = "'6.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.82%  "

# Row 15
$ws.Range("D15").Value = "'2.549.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.73%  "

# Row 16
$ws.Range("D16").Value = "'13.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("D17").Value = "'2.223.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18
$ws.Range("D18").Value = "'0.736"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "

# Row 19
$ws.Range("D19").Value = "'40.031.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0888"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.58%  "

# Row 21
$ws.Range("D21").Value = "'11.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "

# Row 22
$ws.Range("D22").Value = "'5.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.17%  "

# Row 23
$ws.Range("D23").Value = "'65.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "

# Row 24
$ws.Range("D24").Value = "'235.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("E26").Value = "  +2.12%  "

# Row 27
$ws.Range("E27").Value = "  +0.84%  "

# Row 28
$ws.Range("D28").Value = "'23.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.76%  "

# Row 29
$ws.Range("E29").Value = "  -0.96%  "

# Row 30
$ws.Range("D30").Value = "'9.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.84%  "

# Row 31
$ws.Range("D31").Value = "'157.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.73%  "

# Row 32
$ws.Range("D32").Value = "'32.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.52%  "

# Row 33
$ws.Range("E33").Value = "  +0.07%  "

# Row 34
$ws.Range("D34").Value = "'4.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.32%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0715"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.97%  "

# Row 37
$ws.Range("E37").Value = "  -0.64%  "

# Row 38
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.114"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.52%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.102"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.06%  "

# Row 40
$ws.Range("D40").Value = "'1.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.20%  "

# Row 41
$ws.Range("D41").Value = "'15.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "

# Row 42
$ws.Range("D42").Value = "'3.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.37%  "

# Row 43
$ws.Range("D43").Value = "'2.065.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.63%  "

# Row 44
$ws.Range("D44").Value = "'19.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.06%  "

# Row 45
$ws.Range("D45").Value = "'0.0270"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.43%  "

# Row 46
$ws.Range("E46").Value = "  +2.25%  "

# Row 47
$ws.Range("D47").Value = "'2.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.58%  "

# Row 48
$ws.Range("E48").Value = "  -12.38%  "

# Row 49
$ws.Range("D49").Value = "'2.425.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "

# Row 50
$ws.Range("D50").Value = "'1.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.76%  "

# Row 51
$ws.Range("D51").Value = "'1.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.46%  "

